$wb = $excel.ActiveWorkbook
$wsLores = $wb.Worksheets.Item("area_lores")
$wsPopSum = $wb.Worksheets.Item("area_pop_sum")

# ---------------------------------------------------------------------
# New sheet: area_lores_basic - rerun of the density stats for
# area_lores with the "basic" geometry, appended at the end.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLoresBasic = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsLoresBasic.Name = "area_lores_basic"

# The "25%"/"50%"/"75%" index labels must stay literal text and not be
# auto-converted to a percentage number by the smart-entry parser: force
# text format on those cells first, then strip the format back off so no
# stray style lingers on the cells themselves.
$pctLabels = $wsLoresBasic.Range("A6:A8")
$pctLabels.NumberFormat = "@"

$wsLoresBasic.Range("A1").Value = "index"
$wsLoresBasic.Range("B1").Value = "area"

$wsLoresBasic.Range("A2").Value = "count"
$wsLoresBasic.Range("B2").Value = 11

$wsLoresBasic.Range("A3").Value = "mean"
$wsLoresBasic.Range("B3").Value = 9.827399645175809

$wsLoresBasic.Range("A4").Value = "std"
$wsLoresBasic.Range("B4").Value = 5.534037561533252

$wsLoresBasic.Range("A5").Value = "min"
$wsLoresBasic.Range("B5").Value = 2.174119909567626

$wsLoresBasic.Range("A6").Value = "25%"
$wsLoresBasic.Range("B6").Value = 5.217919273144197

$wsLoresBasic.Range("A7").Value = "50%"
$wsLoresBasic.Range("B7").Value = 11.00498254574599

$wsLoresBasic.Range("A8").Value = "75%"
$wsLoresBasic.Range("B8").Value = 14.87710748135346

$wsLoresBasic.Range("A9").Value = "max"
$wsLoresBasic.Range("B9").Value = 16.49826255018992

$pctLabels.ClearFormats()

# Match the bold/border/centered header styling used on the other sheets.
$wsLores.Range("A1:B1").Copy()
$wsLoresBasic.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# New sheet: area_pop_sum_basic - rerun of the density stats for
# area_pop_sum with the "basic" geometry, appended at the end.
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPopSumBasic = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsPopSumBasic.Name = "area_pop_sum_basic"

$wsPopSumBasic.Range("A1").Value = "index"
$wsPopSumBasic.Range("B1").Value = 0

$wsPopSumBasic.Range("A2").Value = "area"
$wsPopSumBasic.Range("B2").Value = 108.1013960969339

$wsPopSumBasic.Range("A3").Value = "population"
$wsPopSumBasic.Range("B3").Value = 190786

$wsPopSumBasic.Range("A4").Value = "density"
$wsPopSumBasic.Range("B4").Value = 1764.880074526728

$wsPopSum.Range("A1:B1").Copy()
$wsPopSumBasic.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
